# ---------------------------------------------------------------------------
# "Cải tiến Tiếp nhận" edit:
#  1. Update existing Sheet1 rows 2 & 3 with refreshed API values and drop a
#     batch of now-empty placeholder cells that the API no longer returns.
#  2. Append two brand-new intake rows (4 & 5) to Sheet1.
#  3. Add a new "Sheet2" summary sheet (FirstName/LastName/InsCardNo/
#     CreateById/InsBenefitRatio/InsBenefitType) for the 4 intake rows.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# ---- Row 2 updates --------------------------------------------------------
$ws1.Range("A2").Value = 1906
$ws1.Range("E2").Value = 46200602599
$ws1.Range("X2").Value = "DN4127460129099"

# Empty placeholder cells the refreshed API payload no longer includes.
$clearRow2 = $ws1.Range("Q2,AE2,AH2,AO2,AX2,AY2,AZ2,BB2,BC2,BD2,BE2,BG2,BH2")
foreach ($area in $clearRow2.Areas) {
    $area.ClearContents()
}

# ---- Row 3 updates --------------------------------------------------------
$ws1.Range("A3").Value = 1907
$ws1.Range("E3").Value = 46200602600
$ws1.Range("AM3").Value = 1
$ws1.Range("AN3").Value = 0

# Row 3 loses its insurance-card block entirely (no matching card this time)
# plus the same batch of empty placeholder cells as row 2.
$clearRow3 = $ws1.Range("Q3,X3,Y3,Z3,AA3,AB3,AC3,AE3,AH3,AO3,AQ3,AX3,AY3,AZ3,BB3,BC3,BD3,BE3,BG3,BH3")
foreach ($area in $clearRow3.Areas) {
    $area.ClearContents()
}

# ---- Row 4 (new intake, with insurance card) ------------------------------
$ws1.Range("A4").Value = 1908
$ws1.Range("B4").Value = "Quách Bảo Hưng"
$ws1.Range("C4").Value = "2000-01-01T00:00:00+07:00"
$ws1.Range("D4").Value = 1
$ws1.Range("E4").Value = 46200602601
$ws1.Range("F4").Value = 356572156
$ws1.Range("G4").Value = "VN"
$ws1.Range("H4").Value = 1
$ws1.Range("I4").Value = "VN"
$ws1.Range("J4").Value = 79
$ws1.Range("K4").Value = 765
$ws1.Range("L4").Value = 26926
$ws1.Range("M4").Value = "5/49 Ntl"
$ws1.Range("N4").Value = 134
$ws1.Range("O4").Value = "Công Ty OPTIMA POWER TOOLS Việt Nam"
$ws1.Range("P4").Value = "Công Ty OPTIMA POWER TOOLS Việt Nam"
$ws1.Range("R4").Value = 1
$ws1.Range("S4").Value = "Quách Bảo Minh Phúc 2"
$ws1.Range("T4").Value = "5/49 Ntl, Phường 07, Quận Bình Thạnh, Thành phố Hồ Chí Minh"
$ws1.Range("U4").Value = 546378305
$ws1.Range("V4").Value = 1
$ws1.Range("W4").Value = "5/49 Ntl"
$ws1.Range("X4").Value = "DN4127460129101"
$ws1.Range("Y4").Value = "BHXH"
$ws1.Range("Z4").Value = "2024-01-01T00:00:00+07:00"
$ws1.Range("AA4").Value = "2025-01-01T00:00:00+07:00"
$ws1.Range("AB4").Value = 3075
$ws1.Range("AC4").Value = 1
$ws1.Range("AD4").Value = 0
$ws1.Range("AF4").Value = $false
$ws1.Range("AG4").Value = $false
$ws1.Range("AI4").Value = 24
$ws1.Range("AJ4").Value = 1
$ws1.Range("AK4").Value = 1
$ws1.Range("AL4").Value = 1
$ws1.Range("AM4").Value = 2
$ws1.Range("AN4").Value = 80
$ws1.Range("AP4").Value = 3839
$ws1.Range("AQ4").Value = "Thẻ BHYT hợp lệ"
$ws1.Range("AR4").Value = 1
$ws1.Range("AS4").Value = 4803
$ws1.Range("AT4").Value = 149
$ws1.Range("AU4").Value = "2024-05-09T09:13:34.2378979+07:00"
$ws1.Range("AV4").Value = 1094172
$ws1.Range("AW4").Value = 560
$ws1.Range("BA4").Value = "||149|09/05/2024 09:13||||||3839||||New|4803|80|||2|1083660|||"
$ws1.Range("BF4").Value = "||1|Normal|CorrectRoute|09/05/2024 09:13|3266971|Quách Bảo Hưng 82|24|Male|01/01/2000 00:00|5/49 Ntl|765|26926|01|VN|134||DN4127389127512|2|80|None|||||||||3839|Open|||||||||||149|09/05/2024 09:13||||||3839||||New|4803|80|||2|1083660|||"

# ---- Row 5 (new intake, no insurance card) --------------------------------
$ws1.Range("A5").Value = 1909
$ws1.Range("B5").Value = "Quách Bảo Hưng"
$ws1.Range("C5").Value = "2000-01-01T00:00:00+07:00"
$ws1.Range("D5").Value = 1
$ws1.Range("E5").Value = 46200602602
$ws1.Range("F5").Value = 356572156
$ws1.Range("G5").Value = "VN"
$ws1.Range("H5").Value = 1
$ws1.Range("I5").Value = "VN"
$ws1.Range("J5").Value = 79
$ws1.Range("K5").Value = 765
$ws1.Range("L5").Value = 26926
$ws1.Range("M5").Value = "5/49 Ntl"
$ws1.Range("N5").Value = 134
$ws1.Range("O5").Value = "Công Ty OPTIMA POWER TOOLS Việt Nam"
$ws1.Range("P5").Value = "Công Ty OPTIMA POWER TOOLS Việt Nam"
$ws1.Range("R5").Value = 1
$ws1.Range("S5").Value = "Quách Bảo Minh Phúc 2"
$ws1.Range("T5").Value = "5/49 Ntl, Phường 07, Quận Bình Thạnh, Thành phố Hồ Chí Minh"
$ws1.Range("U5").Value = 546378305
$ws1.Range("V5").Value = 1
$ws1.Range("W5").Value = "5/49 Ntl"
$ws1.Range("AD5").Value = 0
$ws1.Range("AF5").Value = $false
$ws1.Range("AG5").Value = $false
$ws1.Range("AI5").Value = 24
$ws1.Range("AJ5").Value = 1
$ws1.Range("AK5").Value = 1
$ws1.Range("AL5").Value = 1
$ws1.Range("AM5").Value = 1
$ws1.Range("AN5").Value = 0
$ws1.Range("AP5").Value = 3839
$ws1.Range("AR5").Value = 1
$ws1.Range("AS5").Value = 4803
$ws1.Range("AT5").Value = 149
$ws1.Range("AU5").Value = "2024-05-09T09:13:34.2378979+07:00"
$ws1.Range("AV5").Value = 1094172
$ws1.Range("AW5").Value = 560
$ws1.Range("BA5").Value = "||149|09/05/2024 09:13||||||3839||||New|4803|80|||2|1083660|||"
$ws1.Range("BF5").Value = "||1|Normal|CorrectRoute|09/05/2024 09:13|3266971|Quách Bảo Hưng 82|24|Male|01/01/2000 00:00|5/49 Ntl|765|26926|01|VN|134||DN4127389127512|2|80|None|||||||||3839|Open|||||||||||149|09/05/2024 09:13||||||3839||||New|4803|80|||2|1083660|||"

# ---------------------------------------------------------------------------
# New Sheet2 summary sheet, inserted right after Sheet1.
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Sheet2"

$ws2.Range("A1").Value = "FirstName"
$ws2.Range("B1").Value = "LastName"
$ws2.Range("C1").Value = "InsCardNo"
$ws2.Range("D1").Value = "CreateById"
$ws2.Range("E1").Value = "InsBenefitRatio"
$ws2.Range("F1").Value = "InsBenefitType"

$ws2.Range("A2").Value = 1906
$ws2.Range("B2").Value = "Quách Bảo Hưng"
$ws2.Range("C2").Value = "DN4127460129099"
$ws2.Range("D2").Value = 3839
$ws2.Range("E2").Value = 80
$ws2.Range("F2").Value = 2

$ws2.Range("A3").Value = 1907
$ws2.Range("B3").Value = "Quách Bảo Hưng"
$ws2.Range("D3").Value = 3839
$ws2.Range("E3").Value = 0
$ws2.Range("F3").Value = 1

$ws2.Range("A4").Value = 1908
$ws2.Range("B4").Value = "Quách Bảo Hưng"
$ws2.Range("C4").Value = "DN4127460129101"
$ws2.Range("D4").Value = 3839
$ws2.Range("E4").Value = 80
$ws2.Range("F4").Value = 2

$ws2.Range("A5").Value = 1909
$ws2.Range("B5").Value = "Quách Bảo Hưng"
$ws2.Range("D5").Value = 3839
$ws2.Range("E5").Value = 0
$ws2.Range("F5").Value = 1

# Keep Sheet1 as the active/selected sheet, matching the original workbook view.
$ws1.Activate()
$ws1.Range("A1").Select() | Out-Null

